$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B57: convert from text "3" to a real number 3
$ws.Range("B57").Value = 3

# Add new row 58 with data
$ws.Range("A58").Value = "Ying Tang"
$ws.Range("B58").Value = "'4"
$ws.Range("B58").Style = "Normal"
$ws.Range("C58").Value = "We would like to address "
$ws.Range("D58").Value = "DIS"
$ws.Range("E58").Value = "OTH"
$ws.Range("F58").Value = "f453d9a1-51b4-4aef-ac73-8ca7f3146086"
$ws.Range("G58").Value = "By5SY2gA-_annotated.xlsx"
$ws.Range("H58").Value = "We would like to address the points made individually:"
